$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the en-dash ("\u2013") with a double hyphen ("--") in the three
# affected title strings (rows 3, 7 and 8 of the data table).
$ws.Range("D3").Value = "K25 -- Integrative approach to characterizing gene regulation"
$ws.Range("D7").Value = "Chicago Consortium for Systems Biology -- Core 2: Eukaryotic stress networks"
$ws.Range("D8").Value = "Metaknowledge Network -- Measuring scientific impact"

# Update the selected/active cell from F2 to D9.
$ws.Range("D9").Select()
